$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new boolean column G, importing booleans alongside the existing data
$ws.Range("G1").Value = $true
$ws.Range("G2").Value = $false
$ws.Range("G3").Value = $false
$ws.Range("G4").Value = $false
$ws.Range("G5").Value = $false
$ws.Range("G6").Value = $false
$ws.Range("G7").Value = $false

# Match the selection recorded in the saved file
$ws.Range("K9").Select()
